$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data now starts without a MarketObject: clear the rate-reset block
# (cycleAnchorDateOfRateReset / cycleOfRateReset / marketObjectCodeOfRateReset)
# and zero out the rateSpread column for every annuity row.
for ($row = 2; $row -le 28; $row++) {
    $ws.Range("AA$row").ClearContents()
    $ws.Range("AB$row").ClearContents()
    $ws.Range("AC$row").Value = 0
    $ws.Range("AD$row").ClearContents()
}

# Leave the selection where the author left it after editing the sheet.
$ws.Range("AA29").Select()
